# Onboarding_Email_Handler / Config.xlsx update
# "finalized on #2 and#7"
#
# Adds new Name/Value configuration rows to the "Constants" sheet:
#   - rows 22-26: new SE / Hired / Rejected email settings (subject, body, from address)
#   - rows 27-28: a second Automation-error Email Subject / Body pair
# Also turns the From_Email value (B25) into a mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# ---- long email body texts (kept verbatim, including line breaks) ----

$apos = [char]0x2019   # curly apostrophe (U+2019) used by the source text

$hiredEmailBody = "Dear Candidate`n" + `
  "Thank you for your interest in The Jitu. You have been accepted to proceed to the next stage of the recruitment. We would like to extend our congratulations for making it to this stage..`n" + `
  "Find the attachment for more details.`n" + `
  "Kind Regards,`n" + `
  "Talent Team, TheJitu.`n"

$rejectedEmailBody = "Dear Candidate`n" + `
  "We received an overwhelming number of responses, which makes us feel both humble and proud that so many talented individuals (you included) want to join our team. We know how much effort goes into each application and we appreciate the time taken to contact us. This volume of responses makes for an extremely competitive selection process. Although your profile is impressive, we regret to inform you that we have decided to pursue other candidates for the position at this time.`n" + `
  "However, we are always keen on hearing from talented people and therefore, we strongly encourage you to continue applying for other vacancies advertised on our career website in an area that matches your skillset and experience. Also remember to keep your profile up to date so you can be the first to hear about new job openings.`n" + `
  "Once again, thank you so much for investing your time to make this application. Let${apos}s keep in touch and hopefully we${apos}ll speak again in the near future.`n" + `
  "Kind Regards,`n" + `
  "TheJitu Talent Team."

# ---- column A (Name) values, rows 22-26 : typed first, in row order ----

$ws.Range("A22").Value = "SE_Email_Subject"
$ws.Range("A23").Value = "Hired_Email_Body"
$ws.Range("A24").Value = "Rejected_Email_Subject"
$ws.Range("A25").Value = "From_Email"
$ws.Range("A26").Value = "Hired_Email_Subject"

# ---- column B (Value) values, rows 22-26 : typed next, in row order ----

$ws.Range("B22").Value = "Test Subject"

$ws.Range("B23").Value = $hiredEmailBody
$ws.Range("B23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 14.25

$ws.Range("B24").Value = $rejectedEmailBody
$ws.Range("B24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 14.25

$ws.Range("B25").Value = "emilytiampati@outlook.com"

$ws.Range("B26").Value = "Congratulations You've Been Accepted!"

# ---- rows 27-28 : second Automation-error Subject/Body pair, typed last ----

$ws.Range("A27").Value = "EmailSubject"
$ws.Range("B27").Value = "Automation error!"

$ws.Range("A28").Value = "EmailBody"
$ws.Range("B28").Value = "Hello, "

# ---- turn the From_Email value into a mailto hyperlink ----

$ws.Hyperlinks.Add($ws.Range("B25"), "mailto:emilytiampati@outlook.com")

# ---- restore cursor / selection on the sheet ----

$ws.Range("B33").Select() | Out-Null
